$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Förändrad" column (C) date value changed from 45221 (2023-10-22)
# to 45224 (2023-10-25) for all data rows (rows 2-15).
$ws.Range("C2:C15").Value = 45224
